# Apply the "Add files via upload" update to CORE_holdings.xlsx:
#  - Bump the "as of" date in the confidential disclaimer from 2021-04-06 to 2021-04-08.
#  - Refresh the Weight / Percent Change figures for rows 2-8 on the holdings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected, so it must be unprotected before any cell can be
# changed, and re-protected afterwards with the same password it shipped with.
$ws.Unprotect("D382")

# --- Update the confidential-use disclaimer text (shared string used by A11) ---
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."
# Keep the row's auto height in sync with the (same-length) wrapped text so no
# stray custom-height attribute gets stamped onto the row.
$ws.Rows(11).AutoFit()

# --- Refresh Weight (D) and Percent Change (E) values for rows 2-8 ---
$ws.Range("D2").Value = 0.4942427638340384
$ws.Range("E2").Value = -0.0003490157755128287

$ws.Range("D3").Value = 0.2488981166885851
$ws.Range("E3").Value = 0.0100681077879774

$ws.Range("D4").Value = 0.09802894956437262
$ws.Range("E4").Value = 0.006550768455530598

$ws.Range("D5").Value = 0.1014310908590877
$ws.Range("E5").Value = 0.000289883080490938

$ws.Range("D6").Value = 0.02966778548946358
$ws.Range("E6").Value = 0.001982160555004997

$ws.Range("D7").Value = 0.0277312935644526
$ws.Range("E7").Value = 0.007930954047119076

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = 0.003283744582223358

# Restore sheet protection with the original password.
$ws.Protect("D382")
